# Update the dSF (column F) values to reflect repulled data / recalculated means.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    4  = -4
    6  = 5
    7  = -4
    8  = 3
    9  = 3
    10 = -1
    11 = 1
    13 = 3
    14 = -1
    15 = -7
    16 = 2
    17 = -6
    18 = 2
    19 = -2
    20 = 1
    21 = 1
    22 = 6
    23 = -3
    24 = -6
    25 = 1
    26 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
